$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the R column (year 2022 -> 2023 data) with new figures
$ws.Range("R3").Value = 2023
$ws.Range("R4").Value = 25.6
$ws.Range("R5").Value = 25.6
$ws.Range("R6").Value = 1006
$ws.Range("R7").Value = 971
$ws.Range("R8").Value = 916.6
$ws.Range("R9").Value = 14.3
$ws.Range("R10").Value = 10.7
$ws.Range("R11").Value = 11.3
$ws.Range("R12").Value = 1.1000000000000001
$ws.Range("R13").Value = 13.7

# Update the active selection on the sheet
$ws.Range("R22").Select()
